# Natmi following Dr Hou advice
# Update Ligand/Receptor-expressing cell counts (E, K: 1 -> 3) and
# recompute the dependent NATMI metric columns for rows 2-17.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 57.045267
$ws.Range("H2").Value = 171.135801
$ws.Range("I2").Value = 0.2489699905037019
$ws.Range("J2").Value = 0.2489699905037019
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.7190439999999999
$ws.Range("N2").Value = 2.157132
$ws.Range("O2").Value = 0.03867922735216097
$ws.Range("P2").Value = 0.03867922735216098
$ws.Range("Q2").Value = 41.018056964748
$ws.Range("R2").Value = 369.162512682732
$ws.Range("S2").Value = 0.009629966866558043
$ws.Range("T2").Value = 0.009629966866558044

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 57.045267
$ws.Range("H3").Value = 171.135801
$ws.Range("I3").Value = 0.2489699905037019
$ws.Range("J3").Value = 0.2489699905037019
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 10.81722633333333
$ws.Range("N3").Value = 32.451679
$ws.Range("O3").Value = 0.5818864445941869
$ws.Range("P3").Value = 0.5818864445941871
$ws.Range("Q3").Value = 617.071564384431
$ws.Range("R3").Value = 5553.644079459879
$ws.Range("S3").Value = 0.1448722625848476
$ws.Range("T3").Value = 0.1448722625848476

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 57.045267
$ws.Range("H4").Value = 171.135801
$ws.Range("I4").Value = 0.2489699905037019
$ws.Range("J4").Value = 0.2489699905037019
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 5.999487333333334
$ws.Range("N4").Value = 17.998462
$ws.Range("O4").Value = 0.3227278644455833
$ws.Range("P4").Value = 0.3227278644455833
$ws.Range("Q4").Value = 342.242356793118
$ws.Range("R4").Value = 3080.181211138062
$ws.Range("S4").Value = 0.08034955334629684
$ws.Range("T4").Value = 0.08034955334629686

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 57.045267
$ws.Range("H5").Value = 171.135801
$ws.Range("I5").Value = 0.2489699905037019
$ws.Range("J5").Value = 0.2489699905037019
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 1.054169
$ws.Range("N5").Value = 3.162507
$ws.Range("O5").Value = 0.05670646360806875
$ws.Range("P5").Value = 0.05670646360806875
$ws.Range("Q5").Value = 60.13535206812301
$ws.Range("R5").Value = 541.2181686131071
$ws.Range("S5").Value = 0.01411820770599939
$ws.Range("T5").Value = 0.01411820770599939

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 54.09018966666667
$ws.Range("H6").Value = 162.270569
$ws.Range("I6").Value = 0.2360727666969011
$ws.Range("J6").Value = 0.2360727666969011
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 0.7190439999999999
$ws.Range("N6").Value = 2.157132
$ws.Range("O6").Value = 0.03867922735216097
$ws.Range("P6").Value = 0.03867922735216098
$ws.Range("Q6").Value = 38.89322633867866
$ws.Range("R6").Value = 350.039037048108
$ws.Range("S6").Value = 0.009131112214723094
$ws.Range("T6").Value = 0.009131112214723094

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 54.09018966666667
$ws.Range("H7").Value = 162.270569
$ws.Range("I7").Value = 0.2360727666969011
$ws.Range("J7").Value = 0.2360727666969011
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 10.81722633333333
$ws.Range("N7").Value = 32.451679
$ws.Range("O7").Value = 0.5818864445941869
$ws.Range("P7").Value = 0.5818864445941871
$ws.Range("Q7").Value = 585.1058240372612
$ws.Range("R7").Value = 5265.952416335351
$ws.Range("S7").Value = 0.1373675428787728
$ws.Range("T7").Value = 0.1373675428787728

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 54.09018966666667
$ws.Range("H8").Value = 162.270569
$ws.Range("I8").Value = 0.2360727666969011
$ws.Range("J8").Value = 0.2360727666969011
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 5.999487333333334
$ws.Range("N8").Value = 17.998462
$ws.Range("O8").Value = 0.3227278644455833
$ws.Range("P8").Value = 0.3227278644455833
$ws.Range("Q8").Value = 324.5134077627642
$ws.Range("R8").Value = 2920.620669864878
$ws.Range("S8").Value = 0.07618725984985131
$ws.Range("T8").Value = 0.07618725984985131

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 54.09018966666667
$ws.Range("H9").Value = 162.270569
$ws.Range("I9").Value = 0.2360727666969011
$ws.Range("J9").Value = 0.2360727666969011
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 1.054169
$ws.Range("N9").Value = 3.162507
$ws.Range("O9").Value = 0.05670646360806875
$ws.Range("P9").Value = 0.05670646360806875
$ws.Range("Q9").Value = 57.02020115072034
$ws.Range("R9").Value = 513.181810356483
$ws.Range("S9").Value = 0.01338685175355393
$ws.Range("T9").Value = 0.01338685175355393

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 67.444722
$ws.Range("H10").Value = 202.334166
$ws.Range("I10").Value = 0.2943576685488177
$ws.Range("J10").Value = 0.2943576685488177
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 0.7190439999999999
$ws.Range("N10").Value = 2.157132
$ws.Range("O10").Value = 0.03867922735216097
$ws.Range("P10").Value = 0.03867922735216098
$ws.Range("Q10").Value = 48.49572268576799
$ws.Range("R10").Value = 436.4615041719119
$ws.Range("S10").Value = 0.01138552718465176
$ws.Range("T10").Value = 0.01138552718465177

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 67.444722
$ws.Range("H11").Value = 202.334166
$ws.Range("I11").Value = 0.2943576685488177
$ws.Range("J11").Value = 0.2943576685488177
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 10.81722633333333
$ws.Range("N11").Value = 32.451679
$ws.Range("O11").Value = 0.5818864445941869
$ws.Range("P11").Value = 0.5818864445941871
$ws.Range("Q11").Value = 729.5648228627459
$ws.Range("R11").Value = 6566.083405764713
$ws.Range("S11").Value = 0.1712827371909057
$ws.Range("T11").Value = 0.1712827371909057

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 67.444722
$ws.Range("H12").Value = 202.334166
$ws.Range("I12").Value = 0.2943576685488177
$ws.Range("J12").Value = 0.2943576685488177
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 5.999487333333334
$ws.Range("N12").Value = 17.998462
$ws.Range("O12").Value = 0.3227278644455833
$ws.Range("P12").Value = 0.3227278644455833
$ws.Range("Q12").Value = 404.633755339188
$ws.Range("R12").Value = 3641.703798052692
$ws.Range("S12").Value = 0.09499742175394078
$ws.Range("T12").Value = 0.0949974217539408

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 67.444722
$ws.Range("H13").Value = 202.334166
$ws.Range("I13").Value = 0.2943576685488177
$ws.Range("J13").Value = 0.2943576685488177
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 1.054169
$ws.Range("N13").Value = 3.162507
$ws.Range("O13").Value = 0.05670646360806875
$ws.Range("P13").Value = 0.05670646360806875
$ws.Range("Q13").Value = 71.09813514601801
$ws.Range("R13").Value = 639.8832163141619
$ws.Range("S13").Value = 0.01669198241931949
$ws.Range("T13").Value = 0.0166919824193195

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 50.54489333333333
$ws.Range("H14").Value = 151.63468
$ws.Range("I14").Value = 0.2205995742505793
$ws.Range("J14").Value = 0.2205995742505793
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 0.7190439999999999
$ws.Range("N14").Value = 2.157132
$ws.Range("O14").Value = 0.03867922735216097
$ws.Range("P14").Value = 0.03867922735216098
$ws.Range("Q14").Value = 36.34400228197333
$ws.Range("R14").Value = 327.09602053776
$ws.Range("S14").Value = 0.008532621086228072
$ws.Range("T14").Value = 0.008532621086228074

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 50.54489333333333
$ws.Range("H15").Value = 151.63468
$ws.Range("I15").Value = 0.2205995742505793
$ws.Range("J15").Value = 0.2205995742505793
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 10.81722633333333
$ws.Range("N15").Value = 32.451679
$ws.Range("O15").Value = 0.5818864445941869
$ws.Range("P15").Value = 0.5818864445941871
$ws.Range("Q15").Value = 546.7555511808578
$ws.Range("R15").Value = 4920.79996062772
$ws.Range("S15").Value = 0.128363901939661
$ws.Range("T15").Value = 0.128363901939661

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 50.54489333333333
$ws.Range("H16").Value = 151.63468
$ws.Range("I16").Value = 0.2205995742505793
$ws.Range("J16").Value = 0.2205995742505793
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 5.999487333333334
$ws.Range("N16").Value = 17.998462
$ws.Range("O16").Value = 0.3227278644455833
$ws.Range("P16").Value = 0.3227278644455833
$ws.Range("Q16").Value = 303.2434473180178
$ws.Range("R16").Value = 2729.19102586216
$ws.Range("S16").Value = 0.07119362949549435
$ws.Range("T16").Value = 0.07119362949549436

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 50.54489333333333
$ws.Range("H17").Value = 151.63468
$ws.Range("I17").Value = 0.2205995742505793
$ws.Range("J17").Value = 0.2205995742505793
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 1.054169
$ws.Range("N17").Value = 3.162507
$ws.Range("O17").Value = 0.05670646360806875
$ws.Range("P17").Value = 0.05670646360806875
$ws.Range("Q17").Value = 53.28285966030668
$ws.Range("R17").Value = 479.54573694276
$ws.Range("S17").Value = 0.01250942172919594
$ws.Range("T17").Value = 0.01250942172919594
